# simcon.xlsx - "add more rich info like waste, addvalue, etc."
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix header typo on the Project sheet:
#    CollisionInformationExchnage -> CollisionInformationExchange
# ---------------------------------------------------------------------------
$project = $wb.Worksheets.Item("Project")
$project.Range("G1").Value = "CollisionInformationExchange"

# ---------------------------------------------------------------------------
# 2. Project sheet: update row 2 and append rows 3-5 with more rich info
# ---------------------------------------------------------------------------
$project.Range("A2").Value = 0
$project.Range("G2").Value = 0

$project.Range("A3").Value = 0
$project.Range("B3").Value = 30
$project.Range("C3").Value = 1
$project.Range("D3").Value = 1
$project.Range("E3").Value = 0
$project.Range("F3").Value = 0
$project.Range("G3").Value = 0

$project.Range("A4").Value = 0
$project.Range("B4").Value = 30
$project.Range("C4").Value = 1
$project.Range("D4").Value = 1
$project.Range("E4").Value = 1
$project.Range("F4").Value = 0
$project.Range("G4").Value = 1

$project.Range("A5").Value = 0
$project.Range("B5").Value = 30
$project.Range("C5").Value = 1
$project.Range("D5").Value = 1
$project.Range("E5").Value = 0
$project.Range("F5").Value = 0
$project.Range("G5").Value = 1

$project.Range("G6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Sub sheet: nothing data-wise changes here, just the selection later
# ---------------------------------------------------------------------------
$sub = $wb.Worksheets.Item("Sub")
$sub.Activate() | Out-Null
$sub.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. WorkMethod sheet: add waste (InitialProductionRate all set to 1) and
#    add-value (PerformanceStd set to 0.1) rich info for every work method
# ---------------------------------------------------------------------------
$workMethod = $wb.Worksheets.Item("WorkMethod")
$workMethod.Activate() | Out-Null
For ($r = 2; $r -le 11; $r++) {
    $workMethod.Cells.Item($r, 3).Value = 1
    $workMethod.Cells.Item($r, 5).Value = 0.1
}
$workMethod.Range("E3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. WorkMethodDependency sheet: selection only
# ---------------------------------------------------------------------------
$dependency = $wb.Worksheets.Item("WorkMethodDependency")
$dependency.Activate() | Out-Null
$dependency.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. WorkSpace sheet: selection only
# ---------------------------------------------------------------------------
$workSpace = $wb.Worksheets.Item("WorkSpace")
$workSpace.Activate() | Out-Null
$workSpace.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 7. Task sheet: every task now carries an InitialQty of 5 (the
#    DesignChangeVariation formula column recalculates automatically)
# ---------------------------------------------------------------------------
$task = $wb.Worksheets.Item("Task")
$task.Activate() | Out-Null
For ($r = 2; $r -le 51; $r++) {
    $task.Cells.Item($r, 3).Value = 5
}
$task.Range("C3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 8. Finish back on the Project sheet, which becomes the active tab
# ---------------------------------------------------------------------------
$project.Activate() | Out-Null
$project.Range("G6").Select() | Out-Null
